# Update port breakout v2.1 BOM:
# IC2's part (single inverter buffer -> single Schmitt-Trigger inverter buffer)
# and its MPN (M74VHC1GT04DTT1G -> M74VHC1GT14DTT1G) are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Row 6 holds IC2: B=Part(s) C=Quantity D=Description E=MPN
$ws.Range("E6").Value = "M74VHC1GT14DTT1G"
$ws.Range("D6").Value = "Single Schmitt-Trigger inverter buffer"

# Update the active selection to match the saved state (E6)
$ws.Range("E6").Select()
